$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2: refresh the "CasesTab" row with the new Program-level Cypher query
# (B2) and re-set the still-valid Cases/Arms/etc. count query (C2). A2/D2/E2
# keep their existing text, but we re-assert them for clarity/robustness.
# ---------------------------------------------------------------------------

$casesCountQuery = @'
MATCH (ss:study_subject)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)  
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (ss)<-[:sample_of_study_subject]-(samp)
MATCH (samp)<-[:file_of_sample]-(f)
MATCH (lp)<-[:file_of_laboratory_procedure]-(f)
WHERE p.program_acronym IN ["TAILORx"]
RETURN COUNT(DISTINCT p) AS Programs,
COUNT(DISTINCT s) AS Arms,
COUNT(DISTINCT ss) AS Cases,
COUNT(DISTINCT samp) AS Samples,
COUNT(DISTINCT lp) AS Assays,
COUNT(DISTINCT f) AS Files
'@

$programQuery = @'
MATCH (ss:study_subject)
Match (s:study)
MATCH (ss)<-[:sample_of_study_subject]-(sp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
WITH ss, collect(DISTINCT sp.sample_id) AS samples, 
collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, 
collect(DISTINCT f) AS files
MATCH (s)-[:study_of_program]->(p)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
OPTIONAL MATCH (ss)<-[:program_of_institution]-(p)
OPTIONAL MATCH (p)<-[:of_arm]-(a)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
RETURN DISTINCT
       coalesce (p.program_acronym, '')as `Program Code`,
       coalesce( p.program_id , '')as `Program ID`,
       coalesce (p.program_name , '' )as `Program Name`,
       coalesce(p.start_date, '') as `Start Date`,
       coalesce (p.end_date, '') as `End Date`,
       coalesce(p.pubmed_id, '') as `PubMed ID`,
       count(distinct s) As `Number of Arms`,
       count(distinct ss) as `Associated Cases`
       order By `Program Code`
'@

# Trim the single trailing newline introduced by the here-string literal so
# the stored text matches the source exactly (no trailing blank line).
$casesCountQuery = $casesCountQuery.TrimEnd("`r", "`n")
$programQuery = $programQuery.TrimEnd("`r", "`n")

$ws.Range("A2").Value = "CasesTab"
$ws.Range("B2").Value = $programQuery
$ws.Range("C2").Value = $casesCountQuery
$ws.Range("D2").Value = "TC01_Bento_Filter_Program-TailorX_Neo4jData.xlsx"
$ws.Range("E2").Value = "TC01_Bento_Filter_Program-TailorX_WebData.xlsx"

# Row 2 visually re-wraps taller to fit the longer query now in B2.
$ws.Rows.Item(2).RowHeight = 259.5

# ---------------------------------------------------------------------------
# Row 4 (the old, now-superseded "CasesTab" fragment) is removed entirely.
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).Delete()

# ---------------------------------------------------------------------------
# Column width / sheet formatting tweaks
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 13.02
$ws.Columns.Item(2).ColumnWidth = 74.88
$ws.Columns.Item(4).ColumnWidth = 40.45
$ws.Columns.Item(5).ColumnWidth = 42.88

# ---------------------------------------------------------------------------
# Workbook-level window geometry
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.Left = 31080
$win.Top = 2700
$win.Width = 21600
$win.Height = 11385
